$ErrorActionPreference = "Stop"
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2,1).Value = 1
$ws1.Cells.Item(2,2).Value = '''2024-01-28'
$ws1.Cells.Item(2,3).Value = '苏州.第二届THO 赤维极陵'
$ws1.Cells.Item(2,4).Value = '白塔东路60号(近平江路) 苏州书香府邸平江府'
$ws1.Cells.Item(2,5).Value = '2024.01.28 10:00-01.28 21:00'
$ws1.Cells.Item(2,6).Value = 304
$ws1.Cells.Item(2,7).Value = 65
$ws1.Cells.Item(2,8).Value = 'https://show.bilibili.com/platform/detail.html?id=79002'
$ws1.Cells.Item(2,9).Value = '//i0.hdslb.com/bfs/openplatform/202311/5AgvDWGQ1700817845950.jpeg'

$ws1.Cells.Item(3,1).Value = 2
$ws1.Cells.Item(3,2).Value = '''2024-02-03'
$ws1.Cells.Item(3,3).Value = '【会员购严选】苏州·二次元开放式年会- I COME ACG'
$ws1.Cells.Item(3,4).Value = '金山南路288号木渎影视城F2 苏州广电国际会展中心'
$ws1.Cells.Item(3,5).Value = '2024.02.03 10:00-02.03 20:00'
$ws1.Cells.Item(3,6).Value = 11292
$ws1.Cells.Item(3,7).Value = 25
$ws1.Cells.Item(3,8).Value = 'https://show.bilibili.com/platform/detail.html?id=80426'
$ws1.Cells.Item(3,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/IkyhIHPT1704352086775.jpeg'

$ws1.Cells.Item(4,1).Value = 3
$ws1.Cells.Item(4,2).Value = '''2024-02-03'
$ws1.Cells.Item(4,3).Value = '苏州·TCD国潮动漫游戏嘉年华'
$ws1.Cells.Item(4,4).Value = '苏州大道东688号 苏州国际博览中心'
$ws1.Cells.Item(4,5).Value = '2024.02.03 09:30-02.04 17:00'
$ws1.Cells.Item(4,6).Value = 10584
$ws1.Cells.Item(4,7).Value = 60
$ws1.Cells.Item(4,8).Value = 'https://show.bilibili.com/platform/detail.html?id=80084'
$ws1.Cells.Item(4,9).Value = '//i0.hdslb.com/bfs/openplatform/202401/aDe3s9MS1705479547745.jpeg'

$ws1.Cells.Item(5,1).Value = 4
$ws1.Cells.Item(5,2).Value = '''2024-02-04'
$ws1.Cells.Item(5,3).Value = '苏州·TCD国潮动漫游戏嘉年华吴磊内场'
$ws1.Cells.Item(5,4).Value = '苏州大道东688号 苏州国际博览中心'
$ws1.Cells.Item(5,5).Value = '2024.02.04 09:30-02.04 17:00'
$ws1.Cells.Item(5,6).Value = 598
$ws1.Cells.Item(5,7).Value = '已售罄'
$ws1.Cells.Item(5,8).Value = 'https://show.bilibili.com/platform/detail.html?id=80398'
$ws1.Cells.Item(5,9).Value = '//i1.hdslb.com/bfs/openplatform/202401/bHsHJ3f21704186294427.jpeg'

$ws1.Cells.Item(6,1).Value = 5
$ws1.Cells.Item(6,2).Value = '''2024-02-08'
$ws1.Cells.Item(6,3).Value = '太仓·弇山夜宴'
$ws1.Cells.Item(6,4).Value = '城厢镇县府西街40号公园弄口 弇山园'
$ws1.Cells.Item(6,5).Value = '2024.02.08 17:30-02.24 22:00'
$ws1.Cells.Item(6,6).Value = 1
$ws1.Cells.Item(6,7).Value = 39.9
$ws1.Cells.Item(6,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81215'
$ws1.Cells.Item(6,9).Value = '//i1.hdslb.com/bfs/openplatform/202401/7QA0z2031705908153925.jpeg'

$ws1.Cells.Item(7,1).Value = 6
$ws1.Cells.Item(7,2).Value = '''2024-02-14'
$ws1.Cells.Item(7,3).Value = '常熟·CDW·动漫展02'
$ws1.Cells.Item(7,4).Value = '常熟国际展览中心 国际展览中心'
$ws1.Cells.Item(7,5).Value = '2024.02.14 09:00-02.15 17:30'
$ws1.Cells.Item(7,6).Value = 753
$ws1.Cells.Item(7,7).Value = 55
$ws1.Cells.Item(7,8).Value = 'https://show.bilibili.com/platform/detail.html?id=80504'
$ws1.Cells.Item(7,9).Value = '//i1.hdslb.com/bfs/openplatform/202401/VHHzVjad1704438989848.jpeg'

$ws1.Cells.Item(8,1).Value = 7
$ws1.Cells.Item(8,2).Value = '''2024-02-14'
$ws1.Cells.Item(8,3).Value = '常熟·漫魂动漫游戏展01'
$ws1.Cells.Item(8,4).Value = '虞山北路258号 星程酒店'
$ws1.Cells.Item(8,5).Value = '2024.02.14 09:00-02.14 21:00'
$ws1.Cells.Item(8,6).Value = 106
$ws1.Cells.Item(8,7).Value = 50
$ws1.Cells.Item(8,8).Value = 'https://show.bilibili.com/platform/detail.html?id=80248'
$ws1.Cells.Item(8,9).Value = '//i2.hdslb.com/bfs/openplatform/202312/oPrKUOby1703664065719.jpeg'

$ws1.Cells.Item(9,1).Value = 8
$ws1.Cells.Item(9,2).Value = '''2024-02-14'
$ws1.Cells.Item(9,3).Value = '张家港·META萌元漫展'
$ws1.Cells.Item(9,4).Value = '杨舍镇福新路附近 喜福遇婚庆店'
$ws1.Cells.Item(9,5).Value = '2024.02.14 10:00-02.14 17:00'
$ws1.Cells.Item(9,6).Value = 30
$ws1.Cells.Item(9,7).Value = 20
$ws1.Cells.Item(9,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81189'
$ws1.Cells.Item(9,9).Value = '//i0.hdslb.com/bfs/openplatform/202401/yhLkC15b1705912912966.jpeg'

$ws1.Cells.Item(10,1).Value = 9
$ws1.Cells.Item(10,2).Value = '''2024-02-14'
$ws1.Cells.Item(10,3).Value = '苏州·第一届寒假动漫展宅舞比赛-CF01'
$ws1.Cells.Item(10,4).Value = '润元路润南巷172号,地铁二号线陆慕站东200米,近市旅游换乘中心北100米 斐利酒店'
$ws1.Cells.Item(10,5).Value = '2024.02.14 10:00-02.14 16:00'
$ws1.Cells.Item(10,6).Value = 34
$ws1.Cells.Item(10,7).Value = 49
$ws1.Cells.Item(10,8).Value = 'https://show.bilibili.com/platform/detail.html?id=80528'
$ws1.Cells.Item(10,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/oWbVnOjD1704445446390.jpeg'

$ws1.Cells.Item(11,1).Value = 10
$ws1.Cells.Item(11,2).Value = '''2024-02-16'
$ws1.Cells.Item(11,3).Value = '太仓·龙狮新春动漫节4.0'
$ws1.Cells.Item(11,4).Value = '滨河路126号 凯景世纪大酒店'
$ws1.Cells.Item(11,5).Value = '2024.02.16 08:30-02.16 15:00'
$ws1.Cells.Item(11,6).Value = 30
$ws1.Cells.Item(11,7).Value = 45
$ws1.Cells.Item(11,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81044'
$ws1.Cells.Item(11,9).Value = '//i1.hdslb.com/bfs/openplatform/202401/AMDXVltp1705568031796.jpeg'

$ws1.Cells.Item(12,1).Value = 11
$ws1.Cells.Item(12,2).Value = '''2024-02-16'
$ws1.Cells.Item(12,3).Value = '苏州·Good Jump ACG迎新特别篇X动漫品牌博览会'
$ws1.Cells.Item(12,4).Value = '金山南路288号 广电国际会展中心'
$ws1.Cells.Item(12,5).Value = '2024.02.16 10:00-02.17 17:00'
$ws1.Cells.Item(12,6).Value = 10509
$ws1.Cells.Item(12,7).Value = 60
$ws1.Cells.Item(12,8).Value = 'https://show.bilibili.com/platform/detail.html?id=79303'
$ws1.Cells.Item(12,9).Value = '//i2.hdslb.com/bfs/openplatform/202312/C3P0Encm1701659824998.jpeg'

$ws1.Cells.Item(13,1).Value = 12
$ws1.Cells.Item(13,2).Value = '''2024-02-16'
$ws1.Cells.Item(13,3).Value = '苏州·运动番only专区-Good jump ACG'
$ws1.Cells.Item(13,4).Value = '金山南路288号 广电国际会展中心'
$ws1.Cells.Item(13,5).Value = '2024.02.16 10:00-02.17 17:00'
$ws1.Cells.Item(13,6).Value = 3242
$ws1.Cells.Item(13,7).Value = 25
$ws1.Cells.Item(13,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81435'
$ws1.Cells.Item(13,9).Value = '//i0.hdslb.com/bfs/openplatform/202401/gatL3YjP1706236832019.jpeg'

$ws1.Cells.Item(14,1).Value = 13
$ws1.Cells.Item(14,2).Value = '''2024-02-25'
$ws1.Cells.Item(14,3).Value = '太仓·龙吟动漫游戏展'
$ws1.Cells.Item(14,4).Value = '滨河路128号 凯景世纪大酒店(太仓滨河路店)'
$ws1.Cells.Item(14,5).Value = '2024.02.25 10:00-02.25 17:00'
$ws1.Cells.Item(14,6).Value = 4
$ws1.Cells.Item(14,7).Value = 55
$ws1.Cells.Item(14,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81242'
$ws1.Cells.Item(14,9).Value = '//i1.hdslb.com/bfs/openplatform/202401/swEzpIAK1705915874840.jpeg'

$ws1.Cells.Item(15,1).Value = 14
$ws1.Cells.Item(15,2).Value = '''2024-02-25'
$ws1.Cells.Item(15,3).Value = '苏州·第五届次元鹿角动漫游戏展'
$ws1.Cells.Item(15,4).Value = '清禾路886号 尹山湖大剧院'
$ws1.Cells.Item(15,5).Value = '2024.02.25 10:00-02.25 17:00'
$ws1.Cells.Item(15,6).Value = 2447
$ws1.Cells.Item(15,7).Value = 68
$ws1.Cells.Item(15,8).Value = 'https://show.bilibili.com/platform/detail.html?id=79333'
$ws1.Cells.Item(15,9).Value = '//i1.hdslb.com/bfs/openplatform/202401/tqrMA6qB1704787264871.jpeg'

$ws1.Cells.Item(16,1).Value = 15
$ws1.Cells.Item(16,2).Value = '''2024-03-08'
$ws1.Cells.Item(16,3).Value = '苏州·国风宠物-cosplay展'
$ws1.Cells.Item(16,4).Value = '木渎金山南路288号 苏州国际影视娱乐城'
$ws1.Cells.Item(16,5).Value = '2024.03.08 09:00-03.10 17:30'
$ws1.Cells.Item(16,6).Value = 43
$ws1.Cells.Item(16,7).Value = 65
$ws1.Cells.Item(16,8).Value = 'https://show.bilibili.com/platform/detail.html?id=80635'
$ws1.Cells.Item(16,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/Rfd9PcBN1704781416369.jpeg'

$ws1.Cells.Item(17,1).Value = 16
$ws1.Cells.Item(17,2).Value = '''2024-03-17'
$ws1.Cells.Item(17,3).Value = '苏州·世纪幻想动漫游戏展2.0'
$ws1.Cells.Item(17,4).Value = '清禾路886号 尹山湖大剧院'
$ws1.Cells.Item(17,5).Value = '2024.03.17 10:00-03.17 17:00'
$ws1.Cells.Item(17,6).Value = 19
$ws1.Cells.Item(17,7).Value = 60
$ws1.Cells.Item(17,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81387'
$ws1.Cells.Item(17,9).Value = '//i1.hdslb.com/bfs/openplatform/202401/m0Q2ZB3L1706153205872.jpeg'

$ws1.Cells.Item(18,1).Value = 17
$ws1.Cells.Item(18,2).Value = '''2024-04-13'
$ws1.Cells.Item(18,3).Value = '苏州·绘时国乙1.0-秩序之外'
$ws1.Cells.Item(18,4).Value = '石路步行街永福桥浜15号 银河广场'
$ws1.Cells.Item(18,5).Value = '2024.04.13 13:30-04.13 20:00'
$ws1.Cells.Item(18,6).Value = 92
$ws1.Cells.Item(18,7).Value = 78
$ws1.Cells.Item(18,8).Value = 'https://show.bilibili.com/platform/detail.html?id=80789'
$ws1.Cells.Item(18,9).Value = '//i0.hdslb.com/bfs/openplatform/202401/SjKfDxBh1705041298410.jpeg'

$ws1.Cells.Item(19,1).Value = 18
$ws1.Cells.Item(19,2).Value = '''2024-04-21'
$ws1.Cells.Item(19,3).Value = '苏州·梦幻岛 国乙主题文化展（日夜场） 梦幻岛之约3.0'
$ws1.Cells.Item(19,4).Value = '清禾路888号2号楼3楼 格莱美婚礼宴会中心'
$ws1.Cells.Item(19,5).Value = '2024.04.21 10:00-04.21 21:00'
$ws1.Cells.Item(19,6).Value = 404
$ws1.Cells.Item(19,7).Value = 48.3
$ws1.Cells.Item(19,8).Value = 'https://show.bilibili.com/platform/detail.html?id=78666'
$ws1.Cells.Item(19,9).Value = '//i0.hdslb.com/bfs/openplatform/202312/X0PZ3YhH1703822037665.jpeg'

$ws1.Cells.Item(20,1).Value = 19
$ws1.Cells.Item(20,2).Value = '''2024-05-01'
$ws1.Cells.Item(20,3).Value = '昆山·第十二届理想乡动漫游戏展'
$ws1.Cells.Item(20,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws1.Cells.Item(20,5).Value = '2024.05.01 10:00-05.03 17:00'
$ws1.Cells.Item(20,6).Value = 11078
$ws1.Cells.Item(20,7).Value = 59
$ws1.Cells.Item(20,8).Value = 'https://show.bilibili.com/platform/detail.html?id=77196'
$ws1.Cells.Item(20,9).Value = '//i2.hdslb.com/bfs/openplatform/202310/9xMTQMlg1696736126094.png'

$ws1.Cells.Item(21,1).Value = 20
$ws1.Cells.Item(21,2).Value = '''2024-05-01'
$ws1.Cells.Item(21,3).Value = '苏州·第十七届 I COME ACG  动漫品牌博览会'
$ws1.Cells.Item(21,4).Value = '金山南路288号 广电国际会展中心'
$ws1.Cells.Item(21,5).Value = '2024.05.01 10:00-05.02 17:00'
$ws1.Cells.Item(21,6).Value = 10821
$ws1.Cells.Item(21,7).Value = 65
$ws1.Cells.Item(21,8).Value = 'https://show.bilibili.com/platform/detail.html?id=79789'
$ws1.Cells.Item(21,9).Value = '//i2.hdslb.com/bfs/openplatform/202312/lau3mW031702535438289.jpeg'

$ws1.Cells.Item(22,1).Value = 21
$ws1.Cells.Item(22,2).Value = '''2024-05-02'
$ws1.Cells.Item(22,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾北齐后主签售会'
$ws1.Cells.Item(22,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws1.Cells.Item(22,5).Value = '2024.05.02 14:00-05.02 16:00'
$ws1.Cells.Item(22,6).Value = 8
$ws1.Cells.Item(22,7).Value = 1
$ws1.Cells.Item(22,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81116'
$ws1.Cells.Item(22,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/EubrAneC1705648695005.jpeg'

$ws1.Range("A2").Copy($ws1.Range("A23"))
$ws1.Cells.Item(23,1).Value = 22
$ws1.Cells.Item(23,2).Value = '''2024-05-02'
$ws1.Cells.Item(23,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾啊川签售会'
$ws1.Cells.Item(23,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws1.Cells.Item(23,5).Value = '2024.05.02 14:00-05.02 16:00'
$ws1.Cells.Item(23,6).Value = 17
$ws1.Cells.Item(23,7).Value = 1
$ws1.Cells.Item(23,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81100'
$ws1.Cells.Item(23,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/F24i5GMX1705646667852.jpeg'

$ws1.Range("A2").Copy($ws1.Range("A24"))
$ws1.Cells.Item(24,1).Value = 23
$ws1.Cells.Item(24,2).Value = '''2024-05-02'
$ws1.Cells.Item(24,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾漠小然签售会'
$ws1.Cells.Item(24,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws1.Cells.Item(24,5).Value = '2024.05.02 14:00-05.02 16:00'
$ws1.Cells.Item(24,6).Value = 7
$ws1.Cells.Item(24,7).Value = 1
$ws1.Cells.Item(24,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81119'
$ws1.Cells.Item(24,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/SDnLB1gR1705648838683.jpeg'

$ws1.Range("A2").Copy($ws1.Range("A25"))
$ws1.Cells.Item(25,1).Value = 24
$ws1.Cells.Item(25,2).Value = '''2024-05-02'
$ws1.Cells.Item(25,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾葫芦岛老八签售会'
$ws1.Cells.Item(25,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws1.Cells.Item(25,5).Value = '2024.05.02 14:00-05.02 16:00'
$ws1.Cells.Item(25,6).Value = 6
$ws1.Cells.Item(25,7).Value = 1
$ws1.Cells.Item(25,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81118'
$ws1.Cells.Item(25,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/uHOCneLv1705648779163.jpeg'

$ws1.Range("A2").Copy($ws1.Range("A26"))
$ws1.Cells.Item(26,1).Value = 25
$ws1.Cells.Item(26,2).Value = '''2024-05-03'
$ws1.Cells.Item(26,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾沈辞签售会'
$ws1.Cells.Item(26,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws1.Cells.Item(26,5).Value = '2024.05.03 14:00-05.03 16:00'
$ws1.Cells.Item(26,6).Value = 17
$ws1.Cells.Item(26,7).Value = 1
$ws1.Cells.Item(26,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81120'
$ws1.Cells.Item(26,9).Value = '//i0.hdslb.com/bfs/openplatform/202401/4Pay1rR61705648901961.jpeg'

$ws1.Range("A2").Copy($ws1.Range("A27"))
$ws1.Cells.Item(27,1).Value = 26
$ws1.Cells.Item(27,2).Value = '''2024-05-03'
$ws1.Cells.Item(27,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾矮乐多aliga签售会'
$ws1.Cells.Item(27,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws1.Cells.Item(27,5).Value = '2024.05.03 14:00-05.03 16:00'
$ws1.Cells.Item(27,6).Value = 12
$ws1.Cells.Item(27,7).Value = 1
$ws1.Cells.Item(27,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81114'
$ws1.Cells.Item(27,9).Value = '//i1.hdslb.com/bfs/openplatform/202401/Peub7FOc1705648580577.jpeg'

$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(2).Delete()
$ws2.Rows.Item(2).Delete()

$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2,1).Value = 1
$ws4.Cells.Item(2,2).Value = '''2024-01-28'
$ws4.Cells.Item(2,3).Value = '苏州.第二届THO 赤维极陵'
$ws4.Cells.Item(2,4).Value = '白塔东路60号(近平江路) 苏州书香府邸平江府'
$ws4.Cells.Item(2,5).Value = '2024.01.28 10:00-01.28 21:00'
$ws4.Cells.Item(2,6).Value = 304
$ws4.Cells.Item(2,7).Value = 65
$ws4.Cells.Item(2,8).Value = 'https://show.bilibili.com/platform/detail.html?id=79002'
$ws4.Cells.Item(2,9).Value = '//i0.hdslb.com/bfs/openplatform/202311/5AgvDWGQ1700817845950.jpeg'

$ws4.Cells.Item(3,1).Value = 2
$ws4.Cells.Item(3,2).Value = '''2024-02-03'
$ws4.Cells.Item(3,3).Value = '【会员购严选】苏州·二次元开放式年会- I COME ACG'
$ws4.Cells.Item(3,4).Value = '金山南路288号木渎影视城F2 苏州广电国际会展中心'
$ws4.Cells.Item(3,5).Value = '2024.02.03 10:00-02.03 20:00'
$ws4.Cells.Item(3,6).Value = 11292
$ws4.Cells.Item(3,7).Value = 25
$ws4.Cells.Item(3,8).Value = 'https://show.bilibili.com/platform/detail.html?id=80426'
$ws4.Cells.Item(3,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/IkyhIHPT1704352086775.jpeg'

$ws4.Cells.Item(4,1).Value = 3
$ws4.Cells.Item(4,2).Value = '''2024-02-03'
$ws4.Cells.Item(4,3).Value = '苏州·TCD国潮动漫游戏嘉年华'
$ws4.Cells.Item(4,4).Value = '苏州大道东688号 苏州国际博览中心'
$ws4.Cells.Item(4,5).Value = '2024.02.03 09:30-02.04 17:00'
$ws4.Cells.Item(4,6).Value = 10584
$ws4.Cells.Item(4,7).Value = 60
$ws4.Cells.Item(4,8).Value = 'https://show.bilibili.com/platform/detail.html?id=80084'
$ws4.Cells.Item(4,9).Value = '//i0.hdslb.com/bfs/openplatform/202401/aDe3s9MS1705479547745.jpeg'

$ws4.Cells.Item(5,1).Value = 4
$ws4.Cells.Item(5,2).Value = '''2024-02-04'
$ws4.Cells.Item(5,3).Value = '苏州·TCD国潮动漫游戏嘉年华吴磊内场'
$ws4.Cells.Item(5,4).Value = '苏州大道东688号 苏州国际博览中心'
$ws4.Cells.Item(5,5).Value = '2024.02.04 09:30-02.04 17:00'
$ws4.Cells.Item(5,6).Value = 598
$ws4.Cells.Item(5,7).Value = '已售罄'
$ws4.Cells.Item(5,8).Value = 'https://show.bilibili.com/platform/detail.html?id=80398'
$ws4.Cells.Item(5,9).Value = '//i1.hdslb.com/bfs/openplatform/202401/bHsHJ3f21704186294427.jpeg'

$ws4.Cells.Item(6,1).Value = 5
$ws4.Cells.Item(6,2).Value = '''2024-02-08'
$ws4.Cells.Item(6,3).Value = '太仓·弇山夜宴'
$ws4.Cells.Item(6,4).Value = '城厢镇县府西街40号公园弄口 弇山园'
$ws4.Cells.Item(6,5).Value = '2024.02.08 17:30-02.24 22:00'
$ws4.Cells.Item(6,6).Value = 1
$ws4.Cells.Item(6,7).Value = 39.9
$ws4.Cells.Item(6,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81215'
$ws4.Cells.Item(6,9).Value = '//i1.hdslb.com/bfs/openplatform/202401/7QA0z2031705908153925.jpeg'

$ws4.Cells.Item(7,1).Value = 6
$ws4.Cells.Item(7,2).Value = '''2024-02-14'
$ws4.Cells.Item(7,3).Value = '常熟·CDW·动漫展02'
$ws4.Cells.Item(7,4).Value = '常熟国际展览中心 国际展览中心'
$ws4.Cells.Item(7,5).Value = '2024.02.14 09:00-02.15 17:30'
$ws4.Cells.Item(7,6).Value = 753
$ws4.Cells.Item(7,7).Value = 55
$ws4.Cells.Item(7,8).Value = 'https://show.bilibili.com/platform/detail.html?id=80504'
$ws4.Cells.Item(7,9).Value = '//i1.hdslb.com/bfs/openplatform/202401/VHHzVjad1704438989848.jpeg'

$ws4.Cells.Item(8,1).Value = 7
$ws4.Cells.Item(8,2).Value = '''2024-02-14'
$ws4.Cells.Item(8,3).Value = '常熟·漫魂动漫游戏展01'
$ws4.Cells.Item(8,4).Value = '虞山北路258号 星程酒店'
$ws4.Cells.Item(8,5).Value = '2024.02.14 09:00-02.14 21:00'
$ws4.Cells.Item(8,6).Value = 106
$ws4.Cells.Item(8,7).Value = 50
$ws4.Cells.Item(8,8).Value = 'https://show.bilibili.com/platform/detail.html?id=80248'
$ws4.Cells.Item(8,9).Value = '//i2.hdslb.com/bfs/openplatform/202312/oPrKUOby1703664065719.jpeg'

$ws4.Cells.Item(9,1).Value = 8
$ws4.Cells.Item(9,2).Value = '''2024-02-14'
$ws4.Cells.Item(9,3).Value = '张家港·META萌元漫展'
$ws4.Cells.Item(9,4).Value = '杨舍镇福新路附近 喜福遇婚庆店'
$ws4.Cells.Item(9,5).Value = '2024.02.14 10:00-02.14 17:00'
$ws4.Cells.Item(9,6).Value = 30
$ws4.Cells.Item(9,7).Value = 20
$ws4.Cells.Item(9,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81189'
$ws4.Cells.Item(9,9).Value = '//i0.hdslb.com/bfs/openplatform/202401/yhLkC15b1705912912966.jpeg'

$ws4.Cells.Item(10,1).Value = 9
$ws4.Cells.Item(10,2).Value = '''2024-02-14'
$ws4.Cells.Item(10,3).Value = '苏州·第一届寒假动漫展宅舞比赛-CF01'
$ws4.Cells.Item(10,4).Value = '润元路润南巷172号,地铁二号线陆慕站东200米,近市旅游换乘中心北100米 斐利酒店'
$ws4.Cells.Item(10,5).Value = '2024.02.14 10:00-02.14 16:00'
$ws4.Cells.Item(10,6).Value = 34
$ws4.Cells.Item(10,7).Value = 49
$ws4.Cells.Item(10,8).Value = 'https://show.bilibili.com/platform/detail.html?id=80528'
$ws4.Cells.Item(10,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/oWbVnOjD1704445446390.jpeg'

$ws4.Cells.Item(11,1).Value = 10
$ws4.Cells.Item(11,2).Value = '''2024-02-16'
$ws4.Cells.Item(11,3).Value = '太仓·龙狮新春动漫节4.0'
$ws4.Cells.Item(11,4).Value = '滨河路126号 凯景世纪大酒店'
$ws4.Cells.Item(11,5).Value = '2024.02.16 08:30-02.16 15:00'
$ws4.Cells.Item(11,6).Value = 30
$ws4.Cells.Item(11,7).Value = 45
$ws4.Cells.Item(11,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81044'
$ws4.Cells.Item(11,9).Value = '//i1.hdslb.com/bfs/openplatform/202401/AMDXVltp1705568031796.jpeg'

$ws4.Cells.Item(12,1).Value = 11
$ws4.Cells.Item(12,2).Value = '''2024-02-16'
$ws4.Cells.Item(12,3).Value = '苏州·Good Jump ACG迎新特别篇X动漫品牌博览会'
$ws4.Cells.Item(12,4).Value = '金山南路288号 广电国际会展中心'
$ws4.Cells.Item(12,5).Value = '2024.02.16 10:00-02.17 17:00'
$ws4.Cells.Item(12,6).Value = 10509
$ws4.Cells.Item(12,7).Value = 60
$ws4.Cells.Item(12,8).Value = 'https://show.bilibili.com/platform/detail.html?id=79303'
$ws4.Cells.Item(12,9).Value = '//i2.hdslb.com/bfs/openplatform/202312/C3P0Encm1701659824998.jpeg'

$ws4.Cells.Item(13,1).Value = 12
$ws4.Cells.Item(13,2).Value = '''2024-02-16'
$ws4.Cells.Item(13,3).Value = '苏州·运动番only专区-Good jump ACG'
$ws4.Cells.Item(13,4).Value = '金山南路288号 广电国际会展中心'
$ws4.Cells.Item(13,5).Value = '2024.02.16 10:00-02.17 17:00'
$ws4.Cells.Item(13,6).Value = 3242
$ws4.Cells.Item(13,7).Value = 25
$ws4.Cells.Item(13,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81435'
$ws4.Cells.Item(13,9).Value = '//i0.hdslb.com/bfs/openplatform/202401/gatL3YjP1706236832019.jpeg'

$ws4.Cells.Item(14,1).Value = 13
$ws4.Cells.Item(14,2).Value = '''2024-02-25'
$ws4.Cells.Item(14,3).Value = '太仓·龙吟动漫游戏展'
$ws4.Cells.Item(14,4).Value = '滨河路128号 凯景世纪大酒店(太仓滨河路店)'
$ws4.Cells.Item(14,5).Value = '2024.02.25 10:00-02.25 17:00'
$ws4.Cells.Item(14,6).Value = 4
$ws4.Cells.Item(14,7).Value = 55
$ws4.Cells.Item(14,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81242'
$ws4.Cells.Item(14,9).Value = '//i1.hdslb.com/bfs/openplatform/202401/swEzpIAK1705915874840.jpeg'

$ws4.Cells.Item(15,1).Value = 14
$ws4.Cells.Item(15,2).Value = '''2024-02-25'
$ws4.Cells.Item(15,3).Value = '苏州·第五届次元鹿角动漫游戏展'
$ws4.Cells.Item(15,4).Value = '清禾路886号 尹山湖大剧院'
$ws4.Cells.Item(15,5).Value = '2024.02.25 10:00-02.25 17:00'
$ws4.Cells.Item(15,6).Value = 2447
$ws4.Cells.Item(15,7).Value = 68
$ws4.Cells.Item(15,8).Value = 'https://show.bilibili.com/platform/detail.html?id=79333'
$ws4.Cells.Item(15,9).Value = '//i1.hdslb.com/bfs/openplatform/202401/tqrMA6qB1704787264871.jpeg'

$ws4.Cells.Item(16,1).Value = 15
$ws4.Cells.Item(16,2).Value = '''2024-03-08'
$ws4.Cells.Item(16,3).Value = '苏州·国风宠物-cosplay展'
$ws4.Cells.Item(16,4).Value = '木渎金山南路288号 苏州国际影视娱乐城'
$ws4.Cells.Item(16,5).Value = '2024.03.08 09:00-03.10 17:30'
$ws4.Cells.Item(16,6).Value = 43
$ws4.Cells.Item(16,7).Value = 65
$ws4.Cells.Item(16,8).Value = 'https://show.bilibili.com/platform/detail.html?id=80635'
$ws4.Cells.Item(16,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/Rfd9PcBN1704781416369.jpeg'

$ws4.Cells.Item(17,1).Value = 16
$ws4.Cells.Item(17,2).Value = '''2024-03-17'
$ws4.Cells.Item(17,3).Value = '苏州·世纪幻想动漫游戏展2.0'
$ws4.Cells.Item(17,4).Value = '清禾路886号 尹山湖大剧院'
$ws4.Cells.Item(17,5).Value = '2024.03.17 10:00-03.17 17:00'
$ws4.Cells.Item(17,6).Value = 19
$ws4.Cells.Item(17,7).Value = 60
$ws4.Cells.Item(17,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81387'
$ws4.Cells.Item(17,9).Value = '//i1.hdslb.com/bfs/openplatform/202401/m0Q2ZB3L1706153205872.jpeg'

$ws4.Cells.Item(18,1).Value = 17
$ws4.Cells.Item(18,2).Value = '''2024-04-13'
$ws4.Cells.Item(18,3).Value = '苏州·绘时国乙1.0-秩序之外'
$ws4.Cells.Item(18,4).Value = '石路步行街永福桥浜15号 银河广场'
$ws4.Cells.Item(18,5).Value = '2024.04.13 13:30-04.13 20:00'
$ws4.Cells.Item(18,6).Value = 92
$ws4.Cells.Item(18,7).Value = 78
$ws4.Cells.Item(18,8).Value = 'https://show.bilibili.com/platform/detail.html?id=80789'
$ws4.Cells.Item(18,9).Value = '//i0.hdslb.com/bfs/openplatform/202401/SjKfDxBh1705041298410.jpeg'

$ws4.Cells.Item(19,1).Value = 18
$ws4.Cells.Item(19,2).Value = '''2024-04-21'
$ws4.Cells.Item(19,3).Value = '苏州·梦幻岛 国乙主题文化展（日夜场） 梦幻岛之约3.0'
$ws4.Cells.Item(19,4).Value = '清禾路888号2号楼3楼 格莱美婚礼宴会中心'
$ws4.Cells.Item(19,5).Value = '2024.04.21 10:00-04.21 21:00'
$ws4.Cells.Item(19,6).Value = 404
$ws4.Cells.Item(19,7).Value = 48.3
$ws4.Cells.Item(19,8).Value = 'https://show.bilibili.com/platform/detail.html?id=78666'
$ws4.Cells.Item(19,9).Value = '//i0.hdslb.com/bfs/openplatform/202312/X0PZ3YhH1703822037665.jpeg'

$ws4.Cells.Item(20,1).Value = 19
$ws4.Cells.Item(20,2).Value = '''2024-05-01'
$ws4.Cells.Item(20,3).Value = '昆山·第十二届理想乡动漫游戏展'
$ws4.Cells.Item(20,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws4.Cells.Item(20,5).Value = '2024.05.01 10:00-05.03 17:00'
$ws4.Cells.Item(20,6).Value = 11078
$ws4.Cells.Item(20,7).Value = 59
$ws4.Cells.Item(20,8).Value = 'https://show.bilibili.com/platform/detail.html?id=77196'
$ws4.Cells.Item(20,9).Value = '//i2.hdslb.com/bfs/openplatform/202310/9xMTQMlg1696736126094.png'

$ws4.Cells.Item(21,1).Value = 20
$ws4.Cells.Item(21,2).Value = '''2024-05-01'
$ws4.Cells.Item(21,3).Value = '苏州·第十七届 I COME ACG  动漫品牌博览会'
$ws4.Cells.Item(21,4).Value = '金山南路288号 广电国际会展中心'
$ws4.Cells.Item(21,5).Value = '2024.05.01 10:00-05.02 17:00'
$ws4.Cells.Item(21,6).Value = 10821
$ws4.Cells.Item(21,7).Value = 65
$ws4.Cells.Item(21,8).Value = 'https://show.bilibili.com/platform/detail.html?id=79789'
$ws4.Cells.Item(21,9).Value = '//i2.hdslb.com/bfs/openplatform/202312/lau3mW031702535438289.jpeg'

$ws4.Cells.Item(22,1).Value = 21
$ws4.Cells.Item(22,2).Value = '''2024-05-02'
$ws4.Cells.Item(22,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾北齐后主签售会'
$ws4.Cells.Item(22,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws4.Cells.Item(22,5).Value = '2024.05.02 14:00-05.02 16:00'
$ws4.Cells.Item(22,6).Value = 8
$ws4.Cells.Item(22,7).Value = 1
$ws4.Cells.Item(22,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81116'
$ws4.Cells.Item(22,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/EubrAneC1705648695005.jpeg'

$ws4.Cells.Item(23,1).Value = 22
$ws4.Cells.Item(23,2).Value = '''2024-05-02'
$ws4.Cells.Item(23,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾啊川签售会'
$ws4.Cells.Item(23,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws4.Cells.Item(23,5).Value = '2024.05.02 14:00-05.02 16:00'
$ws4.Cells.Item(23,6).Value = 17
$ws4.Cells.Item(23,7).Value = 1
$ws4.Cells.Item(23,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81100'
$ws4.Cells.Item(23,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/F24i5GMX1705646667852.jpeg'

$ws4.Cells.Item(24,1).Value = 23
$ws4.Cells.Item(24,2).Value = '''2024-05-02'
$ws4.Cells.Item(24,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾漠小然签售会'
$ws4.Cells.Item(24,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws4.Cells.Item(24,5).Value = '2024.05.02 14:00-05.02 16:00'
$ws4.Cells.Item(24,6).Value = 7
$ws4.Cells.Item(24,7).Value = 1
$ws4.Cells.Item(24,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81119'
$ws4.Cells.Item(24,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/SDnLB1gR1705648838683.jpeg'

$ws4.Range("A2").Copy($ws4.Range("A25"))
$ws4.Cells.Item(25,1).Value = 24
$ws4.Cells.Item(25,2).Value = '''2024-05-02'
$ws4.Cells.Item(25,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾葫芦岛老八签售会'
$ws4.Cells.Item(25,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws4.Cells.Item(25,5).Value = '2024.05.02 14:00-05.02 16:00'
$ws4.Cells.Item(25,6).Value = 6
$ws4.Cells.Item(25,7).Value = 1
$ws4.Cells.Item(25,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81118'
$ws4.Cells.Item(25,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/uHOCneLv1705648779163.jpeg'

$ws4.Range("A2").Copy($ws4.Range("A26"))
$ws4.Cells.Item(26,1).Value = 25
$ws4.Cells.Item(26,2).Value = '''2024-05-03'
$ws4.Cells.Item(26,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾沈辞签售会'
$ws4.Cells.Item(26,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws4.Cells.Item(26,5).Value = '2024.05.03 14:00-05.03 16:00'
$ws4.Cells.Item(26,6).Value = 17
$ws4.Cells.Item(26,7).Value = 1
$ws4.Cells.Item(26,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81120'
$ws4.Cells.Item(26,9).Value = '//i0.hdslb.com/bfs/openplatform/202401/4Pay1rR61705648901961.jpeg'

$ws4.Range("A2").Copy($ws4.Range("A27"))
$ws4.Cells.Item(27,1).Value = 26
$ws4.Cells.Item(27,2).Value = '''2024-05-03'
$ws4.Cells.Item(27,3).Value = '昆山·第十二届理想乡动漫游戏展嘉宾矮乐多aliga签售会'
$ws4.Cells.Item(27,4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws4.Cells.Item(27,5).Value = '2024.05.03 14:00-05.03 16:00'
$ws4.Cells.Item(27,6).Value = 12
$ws4.Cells.Item(27,7).Value = 1
$ws4.Cells.Item(27,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81114'
$ws4.Cells.Item(27,9).Value = '//i1.hdslb.com/bfs/openplatform/202401/Peub7FOc1705648580577.jpeg'

Write-Host "edit complete"
